$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 32 with the newly obtained distance-matrix results for the
# 25% subset run (ukb51139_subset.csv, 7003 x 1084).
$ws.Range("A32").Value = "ukb51139_subset.csv"
$ws.Range("B32").Value = "7003 x 1084"
$ws.Range("C32").Value = "all"
$ws.Range("D32").Value = "no events"
$ws.Range("E32").Value = "> 140/80"
$ws.Range("F32").Value = "zscore"
$ws.Range("G32").Value = "median"
$ws.Range("H32").Value = "none"
$ws.Range("I32").Value = 50
$ws.Range("K32").Value = 243
$ws.Range("L32").Value = "-117.4 & -29.7"
$ws.Range("M32").Value = "49.4 & 33.8"
$ws.Range("N32").Value = 16
$ws.Range("O32").Value = 68.9

# The numeric cells in this row (I, K, N, O) pick up the same "normal"
# (non-theme-minor) Calibri 11 font used elsewhere in the sheet once the
# row is populated.
foreach ($addr in @("I32", "K32", "N32", "O32")) {
    $cell = $ws.Range($addr)
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 11
}
